$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows with changed figures.
# NumberFormat is set to Text ("@") before assignment so values such as
# "233.90", "14.60", "0.0930" or "37.257.65" keep their exact original
# textual representation instead of being reinterpreted as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.257.65'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.063.97'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.90'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.618'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.86'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.382'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0762'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.368.35'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.60'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.77'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.778'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.13'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.066.06'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.223.32'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.38'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.45'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0812'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.86'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.42'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.92%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.31'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.12%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.92'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.14%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.117'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.48'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.59'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.86%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.01%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.76'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.83%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.71%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.50%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.476.70'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '95.98'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0930'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.47%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.48%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.16'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.99%  '

# Rows 49-51: coins reshuffled in ranking order - row 49 (was MXToken) becomes
# FraxShare, row 50 (was FraxShare) becomes MXToken, and row 51 (was MultiversX)
# is replaced by RocketPoolETH.
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.17'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.91%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.97'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.34%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.257.37'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.25%  '
